$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update score values (header/login/registration related rows)
$ws.Range("C11").Value = 10   # Web Design
$ws.Range("C15").Value = 3    # Implement Filtering by Status
$ws.Range("C17").Value = 9    # Login Screen

# Fill in previously-empty cells with 0
$ws.Range("C21").Value = 0    # Publish New Ad
$ws.Range("C24").Value = 0    # Implemented Category and Filtering
$ws.Range("C25").Value = 0    # Deactivate/Publish Again Ad
$ws.Range("C26").Value = 0    # Edit Inactive Ads
$ws.Range("C27").Value = 0    # Change and Edit Images
$ws.Range("C28").Value = 0    # Delete Ad
$ws.Range("C29").Value = 0    # Edit User Profile
$ws.Range("C30").Value = 0    # Change Password

# Update the view: scroll position and active selection
$ws.Range("C17").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 13
$window.ScrollColumn = 1
